$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.452.12'
$ws.Range('E2').Value = '  +6.84%  '
$ws.Range('D3').Value = '2.479.30'
$ws.Range('E3').Value = '  +4.77%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '489.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +14.98%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.76%  '
$ws.Range('D9').Value = '2.495.46'
$ws.Range('E9').Value = '  +4.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.79'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0976'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.331'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.02%  '
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '2.911.82'
$ws.Range('E14').Value = '  +4.60%  '
$ws.Range('D15').Value = '56.361.60'
$ws.Range('E15').Value = '  +6.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000136'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.96%  '
$ws.Range('D18').Value = '2.491.65'
$ws.Range('E18').Value = '  +4.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +11.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '318.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  +11.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.410'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.52%  '
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.21%  '
$ws.Range('D28').Value = '2.593.08'
$ws.Range('E28').Value = '  +4.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.62%  '
$ws.Range('D30').Value = '0.0₃0789'
$ws.Range('E30').Value = '  +12.51%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.23'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.20'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.78%  '
$ws.Range('E34').Value = '  +8.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.72'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.19%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.859'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.50'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0560'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.63%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.994'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.610'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.33'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +10.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.82'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +18.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0923'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '258.76'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +19.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0228'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.19'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.86%  '
$ws.Range('D51').Value = '1.889.90'
$ws.Range('E51').Value = '  -1.36%  '
